$d = $word.ActiveDocument

# 1) "Collaborated in a team of 4 to create an ML ..." ->
#    "Collaborated in a team of 4 to publish a ML ..."
$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("create an", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "publish a", 2)

# 2) "Pre-published manuscript in scientific journal" ->
#    "Accepted for publication in scientific journal"
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("Pre-published manuscript", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "Accepted for publication", 2)
